$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new Price-column values below are plain decimal numbers (e.g. 530.93).
# Excel auto-detects such literals as numbers when Value is assigned, but the
# source data keeps these columns as text (note the dotted-thousands values
# like 59.201.31 alongside them), so pin the numeric-looking cells to Text
# format first, matching how the sheet was originally authored.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '59.201.31'
$ws.Range('E2').Value = '  +2.25%  '
$ws.Range('D3').Value = '2.593.36'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '530.93'
$ws.Range('E5').Value = '  +2.86%  '
$ws.Range('D6').Value = '140.24'
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.567'
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('D9').Value = '2.607.28'
$ws.Range('E9').Value = '  +0.94%  '
$ws.Range('D10').Value = '6.45'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('E11').Value = '  +2.75%  '
$ws.Range('D12').Value = '0.334'
$ws.Range('E12').Value = '  +2.12%  '
$ws.Range('E13').Value = '  +3.20%  '
$ws.Range('D14').Value = '3.055.05'
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').Value = '59.133.32'
$ws.Range('E15').Value = '  +2.17%  '
$ws.Range('D16').Value = '20.45'
$ws.Range('E16').Value = '  +1.95%  '
$ws.Range('D17').Value = '0.0000134'
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('D18').Value = '2.598.00'
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('D19').Value = '346.78'
$ws.Range('E19').Value = '  +4.12%  '
$ws.Range('E20').Value = '  +1.30%  '
$ws.Range('D21').Value = '10.11'
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('D22').Value = '6.41'
$ws.Range('E22').Value = '  +1.46%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '67.51'
$ws.Range('E24').Value = '  +2.66%  '
$ws.Range('D25').Value = '0.168'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').Value = '0.406'
$ws.Range('E26').Value = '  +2.09%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '7.16'
$ws.Range('E28').Value = '  +3.20%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = '0.0₃0737'
$ws.Range('E30').Value = '  +1.78%  '
$ws.Range('E31').Value = '  +3.42%  '
$ws.Range('D32').Value = '5.83'
$ws.Range('E32').Value = '  -2.09%  '
$ws.Range('D33').Value = '18.80'
$ws.Range('E33').Value = '  +0.96%  '
$ws.Range('D34').Value = '149.18'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = '3.98'
$ws.Range('E35').Value = '  +1.65%  '
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('D37').Value = '36.86'
$ws.Range('E37').Value = '  +2.09%  '
$ws.Range('E38').Value = '  +4.15%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = '0.831'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('B40').Value = 'SuiNetwork'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('E40').Value = '  +2.03%  '
$ws.Range('D41').Value = '3.53'
$ws.Range('E41').Value = '  +1.83%  '
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').Value = '271.47'
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').Value = '0.0959'
$ws.Range('E46').Value = '  +2.15%  '
$ws.Range('D47').Value = '0.0519'
$ws.Range('E47').Value = '  +1.04%  '
$ws.Range('D48').Value = '18.45'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').Value = '4.65'
$ws.Range('E49').Value = '  +3.18%  '
$ws.Range('D50').Value = '1.951.47'
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('D51').Value = '0.0222'
$ws.Range('E51').Value = '  +1.81%  '
